# Auto-generated Excel COM-interop script
# Applies scheduled-runner market-price updates to the Faerie Profits workbook.
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 1196.5714
$ws.Range("I41").Value = 830.2
$ws.Range("K41").Value = 830.2
$ws.Range("M41").Value = -390.2
$ws.Range("H80").Value = 2617.5833
$ws.Range("I80").Value = 1055.625
$ws.Range("J80").Value = 3398.5625
$ws.Range("K80").Value = 3166.875
$ws.Range("L80").Value = 10195.6875
$ws.Range("M80").Value = -2168.875
$ws.Range("N80").Value = -12191.6875
$ws.Range("H83").Value = 2617.5833
$ws.Range("I83").Value = 1055.625
$ws.Range("J83").Value = 3398.5625
$ws.Range("K83").Value = 9500.625
$ws.Range("L83").Value = 30587.0625
$ws.Range("M83").Value = -4508.625
$ws.Range("N83").Value = -40571.0625
$ws.Range("H106").Value = 1910.88
$ws.Range("I106").Value = 1580.5454
$ws.Range("K106").Value = 1580.5454
$ws.Range("M106").Value = -949.5454
$ws.Range("H137").Value = 2218.077
$ws.Range("I137").Value = 2189.9546
$ws.Range("J137").Value = 2372.75
$ws.Range("K137").Value = 6569.8638
$ws.Range("L137").Value = 7118.25
$ws.Range("M137").Value = -4019.8638
$ws.Range("N137").Value = -12218.25
$ws.Range("I141").Value = 1270.1578
$ws.Range("J141").Value = 7197.4
$ws.Range("K141").Value = 3810.4734
$ws.Range("L141").Value = 21592.2
$ws.Range("M141").Value = 1369.5266
$ws.Range("N141").Value = -31952.2

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6581.264
$ws.Range("I32").Value = 6046.8887
$ws.Range("J32").Value = 9587.125
$ws.Range("K32").Value = 6046.8887
$ws.Range("L32").Value = 9587.125
$ws.Range("M32").Value = -5759.8887
$ws.Range("N32").Value = -10161.125
$ws.Range("H45").Value = 2344.2
$ws.Range("I45").Value = 2089.4614
$ws.Range("K45").Value = 2089.4614
$ws.Range("M45").Value = -1712.4614
$ws.Range("H68").Value = 32777.5
$ws.Range("J68").Value = 45555
$ws.Range("L68").Value = 45555
$ws.Range("N68").Value = -47177
$ws.Range("H71").Value = 32777.5
$ws.Range("J71").Value = 45555
$ws.Range("L71").Value = 136665
$ws.Range("N71").Value = -144777
$ws.Range("H74").Value = 2553.1177
$ws.Range("I74").Value = 1827
$ws.Range("K74").Value = 1827
$ws.Range("M74").Value = -953
$ws.Range("H77").Value = 2553.1177
$ws.Range("I77").Value = 1827
$ws.Range("K77").Value = 9135
$ws.Range("M77").Value = -4767
$ws.Range("H111").Value = 86998.5
$ws.Range("J111").Value = 86998.5
$ws.Range("L111").Value = 86998.5
$ws.Range("N111").Value = -95178.5
$ws.Range("H122").Value = 1506.1818
$ws.Range("I122").Value = 1387.4286
$ws.Range("K122").Value = 4162.2858
$ws.Range("M122").Value = -1712.2858
$ws.Range("H132").Value = 1880.72
$ws.Range("I132").Value = 1376.6086
$ws.Range("K132").Value = 4129.825800000001
$ws.Range("M132").Value = -1599.825800000001

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2165.0278
$ws.Range("I20").Value = 1768.3478
$ws.Range("J20").Value = 2866.8462
$ws.Range("K20").Value = 1768.3478
$ws.Range("L20").Value = 2866.8462
$ws.Range("M20").Value = -1521.3478
$ws.Range("N20").Value = -3360.8462
$ws.Range("H86").Value = 8300.058999999999
$ws.Range("I86").Value = 10874.583
$ws.Range("J86").Value = 2121.2
$ws.Range("K86").Value = 10874.583
$ws.Range("L86").Value = 2121.2
$ws.Range("M86").Value = -9751.583000000001
$ws.Range("N86").Value = -4367.2
$ws.Range("H89").Value = 8300.058999999999
$ws.Range("I89").Value = 10874.583
$ws.Range("J89").Value = 2121.2
$ws.Range("K89").Value = 54372.915
$ws.Range("L89").Value = 10606
$ws.Range("M89").Value = -48756.915
$ws.Range("N89").Value = -21838
$ws.Range("H94").Value = 4120.375
$ws.Range("I94").Value = 3259.6
$ws.Range("K94").Value = 3259.6
$ws.Range("M94").Value = -2808.6
$ws.Range("H107").Value = 11364739
$ws.Range("I107").Value = 13158850
$ws.Range("J107").Value = 2037.3334
$ws.Range("K107").Value = 13158850
$ws.Range("L107").Value = 2037.3334
$ws.Range("M107").Value = -13156930
$ws.Range("N107").Value = -5877.3334

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3181.1924
$ws.Range("I31").Value = 1554.2307
$ws.Range("J31").Value = 4808.154
$ws.Range("K31").Value = 1554.2307
$ws.Range("L31").Value = 4808.154
$ws.Range("M31").Value = -1259.2307
$ws.Range("N31").Value = -5398.154
$ws.Range("H34").Value = 3181.1924
$ws.Range("I34").Value = 1554.2307
$ws.Range("J34").Value = 4808.154
$ws.Range("K34").Value = 1554.2307
$ws.Range("L34").Value = 4808.154
$ws.Range("M34").Value = -1352.2307
$ws.Range("N34").Value = -5212.154
$ws.Range("H69").Value = 13333.333
$ws.Range("I69").Value = 13333.333
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 13333.333
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = -12584.333
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 13333.333
$ws.Range("I72").Value = 13333.333
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 39999.999
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = -36255.999
$ws.Range("N72").ClearContents()

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 102571.82
$ws.Range("J55").Value = 112709
$ws.Range("L55").Value = 338127
$ws.Range("N55").Value = -338481

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 29805.5
$ws.Range("J18").Value = 29805.5
$ws.Range("L18").Value = 29805.5
$ws.Range("N18").Value = -30391.5
$ws.Range("H98").Value = 52133.332
$ws.Range("J98").Value = 52133.332
$ws.Range("L98").Value = 52133.332
$ws.Range("N98").Value = -58123.332
$ws.Range("H104").Value = 82499.5
$ws.Range("J104").Value = 82499.5
$ws.Range("L104").Value = 82499.5
$ws.Range("N104").Value = -89487.5
$ws.Range("H105").Value = 84999.5
$ws.Range("J105").Value = 84999.5
$ws.Range("L105").Value = 84999.5
$ws.Range("N105").Value = -91987.5
$ws.Range("H122").Value = 2045.1904
$ws.Range("I122").Value = 1953.8572
$ws.Range("J122").Value = 2227.8572
$ws.Range("K122").Value = 5861.571599999999
$ws.Range("L122").Value = 6683.571599999999
$ws.Range("M122").Value = -3411.571599999999
$ws.Range("N122").Value = -11583.5716
$ws.Range("H133").Value = 85000
$ws.Range("J133").Value = 85000
$ws.Range("L133").Value = 85000
$ws.Range("N133").Value = -95120

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 27222
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 27222
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 27222
$ws.Range("M20").ClearContents()
$ws.Range("N20").Value = -27674
$ws.Range("H22").Value = 2720.2
$ws.Range("I22").Value = 1179
$ws.Range("K22").Value = 1179
$ws.Range("M22").Value = -884
$ws.Range("H27").Value = 2720.2
$ws.Range("I27").Value = 1179
$ws.Range("K27").Value = 1179
$ws.Range("M27").Value = -1072
$ws.Range("H43").Value = 27222
$ws.Range("J43").Value = 27222
$ws.Range("L43").Value = 27222
$ws.Range("N43").Value = -27608
$ws.Range("H76").Value = 18000
$ws.Range("J76").Value = 18000
$ws.Range("L76").Value = 18000
$ws.Range("N76").Value = -18676
$ws.Range("H79").Value = 18000
$ws.Range("J79").Value = 18000
$ws.Range("L79").Value = 18000
$ws.Range("N79").Value = -20340
$ws.Range("H125").Value = 73183
$ws.Range("J125").Value = 73183
$ws.Range("L125").Value = 73183
$ws.Range("N125").Value = -83023

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 88911.2
$ws.Range("J46").Value = 88911.2
$ws.Range("L46").Value = 88911.2
$ws.Range("N46").Value = -89373.2
$ws.Range("H98").Value = 51590
$ws.Range("J98").Value = 51590
$ws.Range("L98").Value = 51590
$ws.Range("N98").Value = -57580
$ws.Range("H132").Value = 1764.5652
$ws.Range("I132").Value = 1765.0731
$ws.Range("J132").Value = 1760.4
$ws.Range("K132").Value = 5295.219300000001
$ws.Range("L132").Value = 5281.200000000001
$ws.Range("M132").Value = -2765.219300000001
$ws.Range("N132").Value = -10341.2
$ws.Range("H134").Value = 88911.2
$ws.Range("J134").Value = 88911.2
$ws.Range("L134").Value = 266733.6
$ws.Range("N134").Value = -271803.6
